# Expense Details workbook update — add a Travel entry at the top, a
# "koththu" entry below it, and a Travel entry at the bottom (per the
# "open AI was added" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the existing "Food" row (current
# rows 2 and 3), which pushes the four existing expense rows down to
# rows 4-7 and keeps their values untouched.
$ws.Rows("2:3").Insert()

# Row 2: Travel / 1200 / 2025-09-30
$ws.Range("A2").Value = "Travel"
$ws.Range("B2").Value = 1200
$ws.Range("C2").Value = 45930.229537037034

# Row 3: koththu / 2000 / 2025-09-28
$ws.Range("A3").Value = "koththu"
$ws.Range("B3").Value = 2000
$ws.Range("C3").Value = 45928.229537037034

# The inserted rows don't inherit the date number format used by the rest
# of column C, so copy it over from a neighboring date cell.
$ws.Range("C4").Copy()
$ws.Range("C2:C3").PasteSpecial(-4122)  # xlPasteFormats

# Append one more row at the bottom (row 8): Travel / 3000 / 2025-09-01
$ws.Range("A8").Value = "Travel"
$ws.Range("B8").Value = 3000
$ws.Range("C8").Value = 45901.229537037034

$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)  # xlPasteFormats

Write-Output "edit complete"
